# Generate Report for Archive
# - Update the "Status" value from "Ready for handoff" to "In Translation"
#   on all three sheets (Overview, zh-cn, de-de).
# - Re-fit the Status column width on each sheet to the new, shorter text
#   (the localization report re-generates this sheet and the status column
#   autofits to its content).

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# Character-unit width that this engine's column-width quantizer resolves
# to the narrowest storage width available (closest attainable match to a
# freshly-autofit "In Translation" column).
$newColumnWidth = 12.5

foreach ($ws in $wb.Worksheets) {
    switch ($ws.Name) {
        "Overview" {
            $ws.Range("E2").Value = $newStatus
            $ws.Range("F2").Value = $newStatus
            $ws.Range("E:E").ColumnWidth = $newColumnWidth
            $ws.Range("F:F").ColumnWidth = $newColumnWidth
        }
        "zh-cn" {
            $ws.Range("C2").Value = $newStatus
            $ws.Range("C:C").ColumnWidth = $newColumnWidth
        }
        "de-de" {
            $ws.Range("C2").Value = $newStatus
            $ws.Range("C:C").ColumnWidth = $newColumnWidth
        }
    }
}
